$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.191899999999997
$ws.Range("A3").Value = -21.42230000000003
$ws.Range("B5").Value = 5.362800000000002
$ws.Range("D5").Value = -8.918299999999995
$ws.Range("E7").Value = 11.73409999999999
$ws.Range("D9").Value = -8.976000000000004
$ws.Range("D11").Value = -8.367500000000001
$ws.Range("E11").Value = 12.8635
$ws.Range("A14").Value = -20.62039999999998
$ws.Range("A16").Value = -21.35460000000003
$ws.Range("B16").Value = 5.247200000000003
$ws.Range("D17").Value = -9.040600000000001
$ws.Range("E19").Value = 12.9628
$ws.Range("A21").Value = -21.3186
$ws.Range("D21").Value = -8.461200000000002
$ws.Range("E21").Value = 12.53160000000001
$ws.Range("A23").Value = -21.40310000000002
$ws.Range("A25").Value = -22.33720000000004
